$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# The workbook gained a new (unused) named cell style "title_" — a
# bold+underline variant of the existing "title" style — in its style
# table. Register it so the style/font tables match, even though no
# cell on the sheet ends up using it.
$newStyle = $wb.Styles.Add("title_")
$newStyle.Font.Bold = $true
$newStyle.Font.Underline = $true

# Delete row 5 entirely (the "Micro" / "SMEs" / "MSMEs" cells). This also
# drops those now-unused strings from the shared-strings table and shrinks
# the sheet's used range/dimension down to A1:A3.
$ws.Rows(5).Delete()
